# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.369.89'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '3.484.18'
$ws.Range("E3").Value = '  -2.50%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''611.39'
$ws.Range("D6").Value = '''187.83'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("D11").Value = '''52.62'
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("E12").Value = '  -3.76%  '
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D14").Value = '4.044.23'
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("D15").Value = '''603.74'
$ws.Range("E15").Value = '  +3.79%  '
$ws.Range("D16").Value = '69.456.65'
$ws.Range("E16").Value = '  -2.06%  '
$ws.Range("D17").Value = '''18.83'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("D18").Value = '''12.51'
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Value = '3.491.71'
$ws.Range("E19").Value = '  -3.15%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '''0.978'
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").Value = '''17.01'
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").Value = '''105.71'
$ws.Range("E23").Value = '  +12.33%  '
$ws.Range("D24").Value = '''4.70'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").Value = '''5.10'
$ws.Range("E25").Value = '  +4.32%  '
$ws.Range("D26").Value = '''3.00'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").Value = '''10.88'
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").Value = '''9.62'
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("D29").Value = '''33.22'
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("D30").Value = '''6.87'
$ws.Range("E30").Value = '  -4.90%  '
$ws.Range("D31").Value = '''12.47'
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = '''4.02'
$ws.Range("E32").Value = '  +5.62%  '
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("D34").Value = '''63.11'
$ws.Range("D35").Value = '''3.12'
$ws.Range("E35").Value = '  -5.92%  '
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = '3.613.91'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '''3.61'
$ws.Range("E38").Value = '  +4.34%  '
$ws.Range("E39").Value = '  -4.99%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''506.44'
$ws.Range("E40").Value = '  -5.89%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''36.43'
$ws.Range("E41").Value = '  -4.57%  '
$ws.Range("D42").Value = '0.0₃0766'
$ws.Range("E42").Value = '  -7.32%  '
$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D44").Value = '''0.0458'
$ws.Range("E44").Value = '  -3.76%  '
$ws.Range("D45").Value = '''2.88'
$ws.Range("E45").Value = '  -2.70%  '
$ws.Range("D46").Value = '''0.140'
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").Value = '''8.68'
$ws.Range("E49").Value = '  -7.67%  '
$ws.Range("D50").Value = '''131.12'
$ws.Range("E50").Value = '  -2.97%  '
$ws.Range("E51").Value = '  -8.21%  '
